$wb = $excel.ActiveWorkbook

# ============================================================
# Sheet 1: LP1912  (main schedule sheet)
# ============================================================
$ws1 = $wb.Worksheets.Item("LP1912")

# Header: last-updated timestamp + row count
$ws1.Range("A2").Value = "Última actualización: 12:33:21"
$ws1.Range("A3").Value = "Total filas: 162"

# Existing rows whose data was re-ordered by the new scrape run
# (same arrival time / Hora_Llegada, content reshuffled between rows)
# Row 56
$ws1.Range("A56").Value = "08:27:16"
$ws1.Range("C56").Value = "27_EL RETIRO"
$ws1.Range("D56").Value = 50

# Row 57
$ws1.Range("A57").Value = "07:38:39"
$ws1.Range("C57").Value = "14_ABASTO"
$ws1.Range("D57").Value = 99

# Row 86
$ws1.Range("A86").Value = "09:23:23"
$ws1.Range("C86").Value = "27_EL RETIRO"
$ws1.Range("D86").Value = 96

# Row 87
$ws1.Range("A87").Value = "10:50:41"
$ws1.Range("C87").Value = "10_OLMOS"
$ws1.Range("D87").Value = 9

# Row 106
$ws1.Range("A106").Value = "10:05:51"
$ws1.Range("C106").Value = "225_GOMEZ"
$ws1.Range("D106").Value = 107

# Row 108
$ws1.Range("A108").Value = "11:47:17"
$ws1.Range("C108").Value = "23_HERNANDEZ"
$ws1.Range("D108").Value = 5

# Row 135
$ws1.Range("A135").Value = "11:34:59"
$ws1.Range("C135").Value = "23_HERNANDEZ"
$ws1.Range("D135").Value = 62

# Row 136
$ws1.Range("A136").Value = "10:50:41"
$ws1.Range("C136").Value = "27_EL RETIRO"
$ws1.Range("D136").Value = 106

# Row 142
$ws1.Range("A142").Value = "10:50:41"
$ws1.Range("C142").Value = "16_SANTA ANA"
$ws1.Range("D142").Value = 118

# Row 143
$ws1.Range("A143").Value = "11:47:17"
$ws1.Range("C143").Value = "14_ABASTO"
$ws1.Range("D143").Value = 61

# Row 144
$ws1.Range("A144").Value = "11:11:33"
$ws1.Range("C144").Value = "15X38_ABASTO"
$ws1.Range("D144").Value = 97

# Row 146
$ws1.Range("A146").Value = "11:47:17"
$ws1.Range("C146").Value = "11_ETCHEVERRY"
$ws1.Range("D146").Value = 76

# Row 147
$ws1.Range("A147").Value = "11:34:59"
$ws1.Range("C147").Value = "215C_EL PATO"
$ws1.Range("D147").Value = 89

# Row 157
$ws1.Range("A157").Value = "12:11:52"
$ws1.Range("C157").Value = "14_ABASTO"
$ws1.Range("D157").Value = 81

# Row 158
$ws1.Range("A158").Value = "11:34:59"
$ws1.Range("C158").Value = "215A_EL PATO"
$ws1.Range("D158").Value = 118

# Row 162
$ws1.Range("A162").Value = "12:33:21"
$ws1.Range("B162").Value = "13:54"
$ws1.Range("C162").Value = "15_ABASTO"
$ws1.Range("D162").Value = 81

# New rows appended at the end of the sheet (163-167)
$ws1.Range("A163").Value = "12:11:52"
$ws1.Range("B163").Value = "14:01"
$ws1.Range("C163").Value = "10_OLMOS"
$ws1.Range("D163").Value = 110
$ws1.Range("E163").Value = "LP1912"

$ws1.Range("A164").Value = "12:33:21"
$ws1.Range("B164").Value = "14:02"
$ws1.Range("C164").Value = "10_OLMOS"
$ws1.Range("D164").Value = 89
$ws1.Range("E164").Value = "LP1912"

$ws1.Range("A165").Value = "12:33:21"
$ws1.Range("B165").Value = "14:17"
$ws1.Range("C165").Value = "27_EL RETIRO"
$ws1.Range("D165").Value = 104
$ws1.Range("E165").Value = "LP1912"

$ws1.Range("A166").Value = "12:33:21"
$ws1.Range("B166").Value = "14:17"
$ws1.Range("C166").Value = "11_ETCHEVERRY"
$ws1.Range("D166").Value = 104
$ws1.Range("E166").Value = "LP1912"

$ws1.Range("A167").Value = "12:33:21"
$ws1.Range("B167").Value = "14:32"
$ws1.Range("C167").Value = "14X44_ABASTO"
$ws1.Range("D167").Value = 119
$ws1.Range("E167").Value = "LP1912"

# ============================================================
# Sheet 2: LP1912-215 (only the timestamp header changes)
# ============================================================
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 12:33:21"

# ============================================================
# Sheet 3: 6203-6173 (one new row appended)
# ============================================================
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 12:33:21"
$ws3.Range("A3").Value = "Total filas: 23"

# New row 28
$ws3.Range("A28").Value = "12:33:21"
$ws3.Range("B28").Value = "13:57"
$ws3.Range("C28").Value = "215C_LA PLATA"
$ws3.Range("D28").Value = 84
$ws3.Range("E28").Value = "L6203"
